# Insert one new data row above row 58 (pushes old rows 58..154 down to 59..155)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(58).EntireRow.Insert()

# Populate the newly inserted row 58 with its values.
$ws.Range("A58").Value = 7
$ws.Range("B58").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C58").Value = "Ñuble"
$ws.Range("D58").Value = 44477
$ws.Range("E58").Value = 16
$ws.Range("F58").Value = 100112043
$ws.Range("G58").Value = "Pepino ensalada"
$ws.Range("H58").Value = "Sin especificar"
$ws.Range("I58").Value = "Primera"
$ws.Range("J58").Value = 120
$ws.Range("K58").Value = 17000
$ws.Range("L58").Value = 18000
$ws.Range("M58").Value = 17500
$ws.Range("N58").Value = "$/caja 60 unidades"
$ws.Range("O58").Value = "Región de Arica y Parinacota"
$ws.Range("P58").Value = 292
$ws.Range("Q58").Value = 60
$ws.Range("R58").Value = "Hortaliza"

# Make sure the date cell keeps the date number format used by the rest of column D.
$ws.Range("D58").NumberFormat = $ws.Range("D59").NumberFormat()
